# This edit reshuffles the data rows (2-16) of the single worksheet: the
# full row contents (columns D through Q; columns A, B, C, E, F, G, H, R
# are identical across all rows and therefore untouched) of each
# destination row are replaced by the original contents of another row,
# per the mapping below (derived from the target diff).
#
# destination row -> source row (both referring to the ORIGINAL/before state)
#   2  <- 7
#   3  <- 15
#   4  <- 2
#   5  <- 10
#   6  <- 16
#   7  <- 14
#   8  <- 3
#   9  <- 11
#   10 <- 12
#   11 <- 4
#   12 <- 6
#   13 <- 8
#   14 <- 9
#   15 <- 5
#   16 <- 13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 7
    3  = 15
    4  = 2
    5  = 10
    6  = 16
    7  = 14
    8  = 3
    9  = 11
    10 = 12
    11 = 4
    12 = 6
    13 = 8
    14 = 9
    15 = 5
    16 = 13
}

# First, snapshot the current ("before") contents of columns D:Q for every
# data row, since rows will be overwritten in place and we must not read
# already-modified data.
$snapshot = @{}
foreach ($row in 2..16) {
    $snapshot[$row] = $ws.Range("D$row`:Q$row").Value2
}

# Now write each destination row using the snapshotted source row values.
foreach ($destRow in 2..16) {
    $srcRow = $mapping[$destRow]
    $ws.Range("D$destRow`:Q$destRow").Value = $snapshot[$srcRow]
}
